$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Drop the stale "_GoBack" bookmark that used to sit right after the title
#    word "Messung" (it marked the author's last edit position before this
#    revision; Word will re-create it at the new edit point below).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) "... um die MPPT-Ratio zu ermitteln."  ->  "... um den MPP zu ermitteln."
#    (MPPT-Ratio changed to MPP "wehn needed", per the commit message.)
# ---------------------------------------------------------------------------
$oldText = "um die MPPT-Ratio zu ermitteln."
$newText = "um den MPP zu ermitteln."

$d.Content.Find.Execute($oldText, $false, $false, $false, $false, $false, `
                         $true, 1, $false, $newText, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Re-plant "_GoBack" at the point where the new text was typed, i.e. right
#    after "... um den MPP" and before " zu ermitteln." - mirroring Word's
#    own behaviour of stamping the last edit location with this bookmark.
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$sentenceStart = $full.IndexOf($newText)
if ($sentenceStart -ge 0) {
    $markPos = $sentenceStart + "um den MPP".Length
    $markRange = $d.Range($markPos, $markPos)
    $d.Bookmarks.Add("_GoBack", $markRange) | Out-Null
}
